$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Nikroo%Hashemi%NULL%1,            Kathleen%Viveiros%NULL%2,            Walker D.%Redd%NULL%1,            Joyce C.%Zhou%NULL%1,            Thomas R.%McCarty%NULL%1,            Ahmad N.%Bazarbashi%NULL%1,            Kelly E.%Hathorn%NULL%1,            Danny%Wong%NULL%1,            Cheikh%Njie%NULL%1,            Lin%Shen%NULL%1,            Walter W.%Chan%wwchan@bwh.harvard.edu%1]'
$ws.Range("E3").Value = '[Vincent L.%Chen%vichen@med.umich.edu%1,            Fadi%Hawa%NULL%2,            Fadi%Hawa%NULL%0,            Jeffrey A.%Berinstein%NULL%1,            Chanakyaram A.%Reddy%NULL%1,            Ihab%Kassab%NULL%1,            Kevin D.%Platt%NULL%1,            Chia-Yang%Hsu%NULL%1,            Calen A.%Steiner%NULL%1,            Jeremy%Louissaint%NULL%1,            Naresh T.%Gunaratnam%NULL%1,            Pratima%Sharma%pratimas@med.umich.edu%1]'
$ws.Range("E4").Value = '[Carolyn T.%Bramante%NULL%2,            Christopher J.%Tignanelli%NULL%1,            Nirjhar%Dutta%NULL%1,            Emma%Jones%NULL%1,            Leonardo%Tamaritz%NULL%1,            Jeanne%Clark%NULL%1,            Genevieve%Melton-Meaux%NULL%1,            Michael%Usher%NULL%1,            Sayeed%Ikramuddin%NULL%1,             C.%Bramante%null%1,             C. J.% Tignanelli%null%1,             N.% Dutta%null%1,             E.% Jones%null%1,             L.% Tamariz%null%1,             J. M.% Clark%null%1,             M.% Usher%null%1,             G.% Metlon-Meaux%null%1,             S. % Ikramuddin%null%1,        C.%Bramante%null%1,        C. J.% Tignanelli%null%1,        N.% Dutta%null%1,        E.% Jones%null%1,        L.% Tamariz%null%1,        J. M.% Clark%null%1,        M.% Usher%null%1,        G.% Metlon-Meaux%null%1,        S. % Ikramuddin%null%1]'
$ws.Range("E5").Value = '[Kenneth I.%Zheng%NULL%2,            Feng%Gao%NULL%1,            Xiao-Bo%Wang%NULL%2,            Qing-Feng%Sun%NULL%2,            Ke-Hua%Pan%NULL%2,            Ting-Yao%Wang%NULL%2,            Hong-Lei%Ma%NULL%2,            Yong-Ping%Chen%NULL%2,            Wen-Yue%Liu%NULL%1,            Jacob%George%NULL%4,            Ming-Hua%Zheng%NULL%2]'
$ws.Range("E6").Value = '[Rui%Huang%NULL%2,            Li%Zhu%NULL%2,            Li%Zhu%NULL%0,            Jian%Wang%NULL%1,            Leyang%Xue%NULL%1,            Longgen%Liu%NULL%1,            Xuebing%Yan%NULL%1,            Songping%Huang%NULL%1,            Yang%Li%NULL%1,            Xiaomin%Yan%NULL%1,            Biao%Zhang%NULL%1,            Tianmin%Xu%NULL%1,            Chunyang%Li%NULL%1,            Fang%Ji%NULL%1,            Fang%Ming%NULL%1,            Yun%Zhao%NULL%1,            Juan%Cheng%NULL%1,            Yinling%Wang%NULL%1,            Haiyan%Zhao%NULL%1,            Shuqin%Hong%NULL%1,            Kang%Chen%NULL%1,            Xiang‐an%Zhao%NULL%1,            Lei%Zou%NULL%1,            Dawen%Sang%NULL%1,            Huaping%Shao%NULL%1,            Xinying%Guan%NULL%1,            Xiaobing%Chen%NULL%1,            Yuxin%Chen%NULL%1,            Jie%Wei%NULL%1,            Chuanwu%Zhu%zhuchw@126.com%1,            Chao%Wu%dr.wu@nju.edu.cn%1]'
$ws.Range("E7").Value = '[Byrne%Christopher D%coreGivesNoEmail%1,      Chen%Yong-Ping%coreGivesNoEmail%1,      Eslam%Mohammed%coreGivesNoEmail%1,      George%Jacob%coreGivesNoEmail%1,      Mantovani%Alessandro%coreGivesNoEmail%1,      Pan%Ke-Hua%coreGivesNoEmail%1,      Sun%Qing-Feng%coreGivesNoEmail%1,      Targher%Giovanni%coreGivesNoEmail%1,      Wang%Xiao-Bo%coreGivesNoEmail%1,      Yan%Hua-Dong%coreGivesNoEmail%1,      Zheng%Kenneth I%coreGivesNoEmail%1,      Zheng%Ming-Hua%coreGivesNoEmail%1]'
$ws.Range("E8").Value = '[Donghee%Kim%NULL%1,            Nia%Adeniji%NULL%1,            Nyann%Latt%NULL%1,            Sonal%Kumar%NULL%1,            Patricia P.%Bloom%NULL%1,            Elizabeth S.%Aby%NULL%1,            Ponni%Perumalswami%NULL%1,            Marina%Roytman%NULL%1,            Michael%Li%NULL%1,            Alexander S.%Vogel%NULL%1,            Andreea M.%Catana%NULL%1,            Kara%Wegermann%NULL%1,            Rotonya M.%Carr%NULL%1,            Costica%Aloman%NULL%1,            Vincent L.%Chen%NULL%1,            Atoosa%Rabiee%NULL%1,            Brett%Sadowski%NULL%1,            Veronica%Nguyen%NULL%1,            Winston%Dunn%NULL%1,            Kenneth D.%Chavin%NULL%1,            Kali%Zhou%NULL%1,            Blanca%Lizaola-Mayo%NULL%1,            Akshata%Moghe%NULL%1,            José%Debes%NULL%1,            Tzu-Hao%Lee%NULL%1,            Andrea D.%Branch%NULL%1,            Kathleen%Viveiros%NULL%0,            Walter%Chan%NULL%1,            David M.%Chascsa%NULL%1,            Paul%Kwo%NULL%1,            Renumathy%Dhanasekaran%NULL%1]'
$ws.Range("E9").Value = '[Yu-Jie%Zhou%NULL%1,            Kenneth I.%Zheng%NULL%3,            Xiao-Bo%Wang%NULL%0,            Hua-Dong%Yan%NULL%1,            Qing-Feng%Sun%NULL%0,            Ke-Hua%Pan%NULL%0,            Ting-Yao%Wang%NULL%0,            Hong-Lei%Ma%NULL%0,            Yong-Ping%Chen%NULL%0,            Jacob%George%NULL%0,            Ming-Hua%Zheng%NULL%0]'
$ws.Range("E10").Value = '[Yu‐Jie%Zhou%NULL%1,            Kenneth I.%Zheng%NULL%0,            Kenneth I.%Zheng%NULL%0,            Xiao‐Bo%Wang%NULL%4,            Xiao‐Bo%Wang%NULL%0,            Qing‐Feng%Sun%NULL%2,            Ke‐Hua%Pan%NULL%2,            Ting‐Yao%Wang%NULL%2,            Hong‐Lei%Ma%NULL%1,            Yong‐Ping%Chen%NULL%2,            Jacob%George%jacob.george@sydney.edu.au%0,            Ming‐Hua%Zheng%zhengmh@wmu.edu.cn%2]'
$ws.Range("E11").Value = '[Feng%Gao%NULL%1,            Kenneth I%Zheng%NULL%1,            Xiao‐Bo%Wang%NULL%0,            Xiao‐Bo%Wang%NULL%0,            Hua‐Dong%Yan%NULL%1,            Qing‐Feng%Sun%NULL%0,            Ke‐Hua%Pan%NULL%0,            Ting‐Yao%Wang%NULL%0,            Yong‐Ping%Chen%NULL%0,            Jacob%George%jacob.george@sydney.edu.au%0,            Ming‐Hua%Zheng%zhengmh@wmu.edu.cn%0]'
$ws.Range("E12").Value = '[Ivan%Lopez-Mendez%NULL%1,            Jorge%Aquino-Matus%NULL%1,            Sofia Murua-Beltrán%Gall%NULL%1,            Jose D.%Prieto-Nava%NULL%1,            Eva%Juarez-Hernandez%NULL%1,            Misael%Uribe%NULL%1,            Graciela%Castro-Narro%NULL%1]'
$ws.Range("E13").Value = '[Roberta%Forlano%NULL%1,            Benjamin H.%Mullish%NULL%2,            Benjamin H.%Mullish%NULL%0,            Sujit K.%Mukherjee%NULL%1,            Rooshi%Nathwani%NULL%1,            Cristopher%Harlow%NULL%1,            Peter%Crook%NULL%1,            Rebekah%Judge%NULL%1,            Anet%Soubieres%NULL%2,            Anet%Soubieres%NULL%0,            Paul%Middleton%NULL%1,            Anna%Daunt%NULL%1,            Pablo%Perez-Guzman%NULL%1,            Nowlan%Selvapatt%NULL%1,            Maud%Lemoine%NULL%1,            Ameet%Dhar%NULL%1,            Mark R.%Thursz%NULL%1,            Shevanthi%Nayagam%NULL%1,            Pinelopi%Manousou%NULL%1,            Wan-Long%Chuang%NULL%2,            Wan-Long%Chuang%NULL%0,            NULL%NULL%NULL%3,            NULL%NULL%NULL%0]'
$ws.Range("E14").Value = '[Mahmud%Mahamid%NULL%1,            Wiliam%Nseir%NULL%1,            Tawfik%Khoury%NULL%1,            Baker%Mahamid%NULL%1,            Adi%Nubania%NULL%1,            Kamal%Sub-Laban%NULL%1,            Joel%Schifter%NULL%1,            Amir%Mari%NULL%1,            Wisam%Sbeit%NULL%1,            Eran%Goldin%NULL%1]'
